$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells I1 (I0) and J1 (IF) ----------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the existing header style (same as H1 / "IP") instead of inventing a
# brand-new style entry: copy H1's formatting onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data columns I2:J86 ------------------------------------------------------
$iVals = @(6,8,6,7,8,7,8,7,10,7,7,8,8,7,8,8,8,8,8,8,6,8,8,6,7,7,7,6,7,7,7,7,7,7,7,8,7,7,7,8,7,9,7,7,8,8,7,7,7,7,8,8,7,8,7,9,7,7,8,7,7,8,6,8,7,7,7,6,8,7,7,7,7,7,7,9,6,7,7,6,7,6,7,5,6)
$jVals = @(7,8,7,8,8,7,8,7,11,7,7,8,8,7,8,8,8,8,8,8,7,8,8,7,7,7,7,7,7,7,7,7,7,7,7,8,7,7,7,8,7,9,7,7,8,8,7,7,8,7,8,8,7,8,7,9,7,7,8,7,7,8,7,8,7,7,7,7,8,7,7,8,7,7,7,9,6,8,7,6,7,6,7,5,6)

for ($k = 0; $k -lt $iVals.Length; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$k]
    $ws.Cells.Item($row, 10).Value = $jVals[$k]
}
